$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 2375.553533848452
$ws.Range("C3").Value = 2343.776963035138
$ws.Range("C4").Value = 2375.553533848452
$ws.Range("C5").Value = 2375.553533848452
$ws.Range("C6").Value = 2348.940742550856
$ws.Range("C7").Value = 2348.940742550856
$ws.Range("C8").Value = 2042.996992853785
$ws.Range("C9").Value = 2054.322229650993
$ws.Range("C10").Value = 2029.362860345671
$ws.Range("C11").Value = 2039.511913560559
$ws.Range("C12").Value = 1978.054592518947
$ws.Range("C13").Value = 1942.811553016351
$ws.Range("C14").Value = 20060.98989507073
$ws.Range("C15").Value = 20172.47624064033
$ws.Range("C16").Value = 20172.47624064033
$ws.Range("C17").Value = 19663.49715172883
$ws.Range("C18").Value = 20252.92847467912
$ws.Range("C19").Value = 20442.06886548464
$ws.Range("C20").Value = 1556.628165139704
$ws.Range("C21").Value = 1556.628165139704
$ws.Range("C22").Value = 1497.392386900223
$ws.Range("C23").Value = 1433.340338412169
$ws.Range("C24").Value = 1400.948862870578
$ws.Range("C25").Value = 1456.875815862063
$ws.Range("C26").Value = 2472.622910278986
$ws.Range("C27").Value = 2223.999680159021
$ws.Range("C28").Value = 2223.999680159021
$ws.Range("C29").Value = 2192.917481138555
$ws.Range("C30").Value = 2170.395154017165
$ws.Range("C31").Value = 2180.639637307194
$ws.Range("C32").Value = 1269.752317301101
$ws.Range("C33").Value = 1264.57568970012
$ws.Range("C34").Value = 1264.57568970012
$ws.Range("C35").Value = 1179.823258639975
$ws.Range("C36").Value = 1226.182278134155
$ws.Range("C37").Value = 1212.954307496976
$ws.Range("C38").Value = 1195.469426699615
$ws.Range("C39").Value = 1149.435295122411
$ws.Range("C40").Value = 1164.417946995289
$ws.Range("C41").Value = 1121.200725318295
$ws.Range("C42").Value = 1121.341595103268
$ws.Range("C43").Value = 1125.444185667789
$ws.Range("C44").Value = 921.1366307298617
$ws.Range("C45").Value = 921.1366307298617
$ws.Range("C46").Value = 921.1366307298617
$ws.Range("C47").Value = 899.664551057899
$ws.Range("C48").Value = 942.9843646082304
$ws.Range("C49").Value = 942.9843646082304
$ws.Range("C50").Value = 1326.104443508854
$ws.Range("C51").Value = 1326.104443508854
$ws.Range("C52").Value = 1310.176744372151
$ws.Range("C53").Value = 1304.174808347777
$ws.Range("C54").Value = 1274.393002091884
$ws.Range("C55").Value = 1278.495592656406
$ws.Range("C56").Value = 1110.533486579043
$ws.Range("C57").Value = 1055.029172650482
$ws.Range("C58").Value = 1185.823539096346
$ws.Range("C59").Value = 1166.590835918455
$ws.Range("C60").Value = 1166.590835918455
$ws.Range("C61").Value = 1170.693426482977
$ws.Range("C62").Value = 1549.711507040808
$ws.Range("C63").Value = 1625.096184762227
$ws.Range("C64").Value = 1625.096184762227
$ws.Range("C65").Value = 1444.442677860478
$ws.Range("C66").Value = 1348.542213620639
$ws.Range("C67").Value = 1326.971640725344
$ws.Range("C68").Value = 1128.813106177955
$ws.Range("C69").Value = 1125.374733954056
$ws.Range("C70").Value = 1102.20540222471
$ws.Range("C71").Value = 1125.374733954056
$ws.Range("C72").Value = 1048.742565567924
$ws.Range("C73").Value = 1048.742565567924
$ws.Range("C74").Value = 4008.622999057393
$ws.Range("C75").Value = 3586.323824722017
$ws.Range("C76").Value = 3879.233787054284
$ws.Range("C77").Value = 3168.300350252005
$ws.Range("C78").Value = 3193.629945289124
$ws.Range("C79").Value = 3463.454956561318
$ws.Range("C80").Value = 265.2727661758688
$ws.Range("C81").Value = 237.0015010780877
$ws.Range("C82").Value = 353.0525387537038
$ws.Range("C83").Value = 353.0525387537038
$ws.Range("C84").Value = 328.5000562691307
$ws.Range("C85").Value = 328.5000562691307
$ws.Range("C86").Value = 2679.102706171934
$ws.Range("C87").Value = 2733.933373082208
$ws.Range("C88").Value = 2910.05686265333
$ws.Range("C89").Value = 2791.22993329431
$ws.Range("C90").Value = 2829.011456282065
$ws.Range("C91").Value = 2892.639761637603
$ws.Range("C92").Value = 13871.39536009172
$ws.Range("C93").Value = 14275.03494417892
$ws.Range("C94").Value = 14168.98056215686
$ws.Range("C95").Value = 14179.64168570854
$ws.Range("C96").Value = 14391.8918163149
$ws.Range("C97").Value = 14281.73484372832
$ws.Range("C98").Value = 5227.908173271192
$ws.Range("C99").Value = 5227.908173271192
$ws.Range("C100").Value = 5227.908173271192
$ws.Range("C101").Value = 5215.309203997353
$ws.Range("C102").Value = 5187.930582265887
$ws.Range("C103").Value = 5190.885592995852
$ws.Range("C104").Value = 907.7141737027347
$ws.Range("C105").Value = 822.0709476746173
$ws.Range("C106").Value = 822.0709476746173
$ws.Range("C107").Value = 822.0709476746173
$ws.Range("C108").Value = 822.0709476746173
$ws.Range("C109").Value = 822.0709476746173
$ws.Range("C110").Value = 426.484944206429
$ws.Range("C111").Value = 426.484944206429
$ws.Range("C112").Value = 426.484944206429
$ws.Range("C113").Value = 419.1502344715777
$ws.Range("C114").Value = 404.0814408338649
$ws.Range("C115").Value = 419.1502344715777
$ws.Range("C116").Value = 348.9146051669941
$ws.Range("C117").Value = 350.9610384721077
$ws.Range("C118").Value = 383.9095616100801
$ws.Range("C119").Value = 383.9095616100801
$ws.Range("C120").Value = 411.5208861088303
$ws.Range("C121").Value = 538.8805350539859
$ws.Range("C122").Value = 2716.030410733037
$ws.Range("C123").Value = 2805.990132450204
$ws.Range("C124").Value = 2772.803572747229
$ws.Range("C125").Value = 2657.335120812977
$ws.Range("C126").Value = 2657.335120812977
$ws.Range("C127").Value = 2810.719920762579
$ws.Range("C128").Value = 158.3982980600248
$ws.Range("C129").Value = 151.9805821906199
$ws.Range("C130").Value = 153.8823344314395
$ws.Range("C131").Value = 153.4672872230555
$ws.Range("C132").Value = 151.9805821906199
$ws.Range("C133").Value = 151.9805821906199
$ws.Range("C134").Value = 1001.735284150767
$ws.Range("C135").Value = 1007.099994460095
$ws.Range("C136").Value = 1302.755614110094
$ws.Range("C137").Value = 1315.391497500657
$ws.Range("C138").Value = 1403.622790975014
$ws.Range("C139").Value = 1394.037641173902
